$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 10-12, column A: give the label cells the same "mtitleStyle"
# formatting already used by the header row (row 9 / row 15), i.e. s="4"
# in the OOXML. Copy/PasteSpecial (formats only) reuses the existing style
# instead of registering a brand-new (duplicate) style record.
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- C11: "-3" -> "-1" (stays a text value, not a number).
# A leading apostrophe forces Excel to store it as text instead of
# auto-converting the numeric-looking string to a number.
$ws.Range("C11").Value = "'-1"
# Re-apply the original cell formatting (style s="6", same as before),
# since writing the text value above can otherwise nudge the cell onto a
# differently-flagged (quote-prefixed) style record.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- C12: -21 -> -7 (numeric).
$ws.Range("C12").Value = -7

# --- E12: "79/140" -> "93/140" (text).
$ws.Range("E12").Value = "93/140"
